$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("J5").Value = 1.02
$ws.Range("K5").Value = 19

# Row 15
$ws.Range("I15").Value = 3.25
$ws.Range("X15").Value = 19
$ws.Range("AD15").Value = 301
$ws.Range("AE15").Value = 9
$ws.Range("AI15").Value = 29

# Row 16
$ws.Range("J16").Value = 1.08
$ws.Range("K16").Value = 8

# Row 17
$ws.Range("G17").Value = 2.88
$ws.Range("H17").Value = 3.3
$ws.Range("I17").Value = 2.3
$ws.Range("J17").Value = 1.05
$ws.Range("K17").Value = 11
$ws.Range("X17").Value = 23
$ws.Range("AA17").Value = 6.5
$ws.Range("AE17").Value = 8.5
$ws.Range("AF17").Value = 12
$ws.Range("AG17").Value = 9.5
$ws.Range("AI17").Value = 19

# Row 18
$ws.Range("H18").Value = 3.9
$ws.Range("N18").Value = 1.9
$ws.Range("O18").Value = 1.9
$ws.Range("U18").Value = 7
$ws.Range("AA18").Value = 7.5
$ws.Range("AC18").Value = 51
$ws.Range("AD18").Value = 351
$ws.Range("AG18").Value = 19

# Row 19
$ws.Range("N19").Value = 2.08
$ws.Range("O19").Value = 1.73

# Row 26
$ws.Range("G26").Value = 1.14
$ws.Range("T26").Value = 6
$ws.Range("V26").Value = 11
$ws.Range("AE26").Value = 34
$ws.Range("AF26").Value = 81
$ws.Range("AG26").Value = 51
$ws.Range("AH26").Value = 401
$ws.Range("AI26").Value = 201
$ws.Range("AJ26").Value = 151

# Row 27
$ws.Range("H27").Value = 3.1
$ws.Range("I27").Value = 3.4
$ws.Range("J27").Value = 1.08
$ws.Range("K27").Value = 6.7
$ws.Range("L27").Value = 1.34
$ws.Range("M27").Value = 3
$ws.Range("N27").Value = 2
$ws.Range("O27").Value = 1.72
$ws.Range("P27").Value = 1.47
$ws.Range("Q27").Value = 2.5
$ws.Range("S27").Value = 1.95
$ws.Range("T27").Value = 7.4
$ws.Range("U27").Value = 10.5
$ws.Range("X27").Value = 17.5
$ws.Range("Y27").Value = 28
$ws.Range("Z27").Value = 6.7
$ws.Range("AA27").Value = 6
$ws.Range("AH27").Value = 50
$ws.Range("AJ27").Value = 40

# Row 29
$ws.Range("G29").Value = 2.72
$ws.Range("H29").Value = 3.1
$ws.Range("I29").Value = 2.57
$ws.Range("K29").Value = 6.4
$ws.Range("L29").Value = 1.38
$ws.Range("M29").Value = 2.8
$ws.Range("N29").Value = 2.12
$ws.Range("P29").Value = 1.5
$ws.Range("Q29").Value = 2.42
$ws.Range("R29").Value = 1.83
$ws.Range("S29").Value = 1.87
$ws.Range("T29").Value = 7.8
$ws.Range("U29").Value = 13
$ws.Range("V29").Value = 10
$ws.Range("W29").Value = 32
$ws.Range("X29").Value = 25
$ws.Range("Y29").Value = 35
$ws.Range("Z29").Value = 6.4
$ws.Range("AA29").Value = 5.9
$ws.Range("AB29").Value = 14.5
$ws.Range("AC29").Value = 75
$ws.Range("AE29").Value = 7.5
$ws.Range("AF29").Value = 12
$ws.Range("AG29").Value = 9.75
$ws.Range("AH29").Value = 28
$ws.Range("AI29").Value = 23
$ws.Range("AJ29").Value = 35

# Row 30
$ws.Range("G30").Value = 2.6
$ws.Range("H30").Value = 2.92
$ws.Range("I30").Value = 2.82
$ws.Range("R30").Value = 1.91
$ws.Range("S30").Value = 1.8
$ws.Range("T30").Value = 6.7
$ws.Range("U30").Value = 11.75
$ws.Range("V30").Value = 10
$ws.Range("W30").Value = 29
$ws.Range("X30").Value = 26
$ws.Range("Y30").Value = 40
$ws.Range("AB30").Value = 15
$ws.Range("AC30").Value = 80
$ws.Range("AE30").Value = 7.6
$ws.Range("AF30").Value = 13.5
$ws.Range("AG30").Value = 10.25
$ws.Range("AH30").Value = 35
$ws.Range("AI30").Value = 26
$ws.Range("AJ30").Value = 37

# Row 32
$ws.Range("J32").Value = 1.07
$ws.Range("K32").Value = 9
$ws.Range("Z32").Value = 9
$ws.Range("AD32").Value = 351

# Row 34
$ws.Range("G34").Value = 10.75
$ws.Range("T34").Value = 23
$ws.Range("U34").Value = 65
$ws.Range("V34").Value = 29
$ws.Range("W34").Value = 250
$ws.Range("X34").Value = 110
$ws.Range("AB34").Value = 22

# Row 36
$ws.Range("G36").Value = 1.95
$ws.Range("H36").Value = 3.3
$ws.Range("I36").Value = 3.65
$ws.Range("L36").Value = 1.4
$ws.Range("M36").Value = 2.52
$ws.Range("N36").Value = 2.15
$ws.Range("O36").Value = 1.55
$ws.Range("P36").Value = 1.47
$ws.Range("Q36").Value = 2.32
$ws.Range("R36").Value = 1.98
$ws.Range("S36").Value = 1.65
$ws.Range("T36").Value = 5.9
$ws.Range("U36").Value = 8.25
$ws.Range("V36").Value = 9
$ws.Range("W36").Value = 16.5
$ws.Range("X36").Value = 18
$ws.Range("Y36").Value = 37
$ws.Range("Z36").Value = 7.8
$ws.Range("AA36").Value = 6.5
$ws.Range("AB36").Value = 18.5
$ws.Range("AC36").Value = 110
$ws.Range("AD36").Value = 1000
$ws.Range("AE36").Value = 8.75
$ws.Range("AF36").Value = 17.5
$ws.Range("AG36").Value = 13
$ws.Range("AH36").Value = 50
$ws.Range("AI36").Value = 40

# Row 37
$ws.Range("G37").Value = 2.25
$ws.Range("I37").Value = 3
$ws.Range("L37").Value = 1.42
$ws.Range("M37").Value = 2.45
$ws.Range("N37").Value = 2.22
$ws.Range("O37").Value = 1.52
$ws.Range("P37").Value = 1.5
$ws.Range("Q37").Value = 2.25
$ws.Range("R37").Value = 1.98
$ws.Range("S37").Value = 1.65
$ws.Range("T37").Value = 6.3
$ws.Range("U37").Value = 9.75
$ws.Range("V37").Value = 9.5
$ws.Range("W37").Value = 22
$ws.Range("X37").Value = 22
$ws.Range("Z37").Value = 7.4
$ws.Range("AB37").Value = 18
$ws.Range("AC37").Value = 110
$ws.Range("AE37").Value = 7.5
$ws.Range("AF37").Value = 14
$ws.Range("AG37").Value = 11.5
$ws.Range("AH37").Value = 37
$ws.Range("AI37").Value = 32
$ws.Range("AJ37").Value = 50

# Row 39
$ws.Range("N39").Value = 1.93
$ws.Range("O39").Value = 1.88

# Row 44
$ws.Range("G44").Value = 2.5
$ws.Range("I44").Value = 2.75
$ws.Range("T44").Value = 8.5
$ws.Range("U44").Value = 12
$ws.Range("W44").Value = 23
$ws.Range("AE44").Value = 9.5
$ws.Range("AH44").Value = 29

# Row 53
$ws.Range("J53").Value = 1.05
$ws.Range("K53").Value = 11
$ws.Range("N53").Value = 1.85
$ws.Range("O53").Value = 2

# Row 55
$ws.Range("N55").Value = 1.73
$ws.Range("O55").Value = 2.08

# Row 61
$ws.Range("N61").Value = 1.73
$ws.Range("O61").Value = 2.08

# Row 63
$ws.Range("N63").Value = 1.92
$ws.Range("O63").Value = 1.82

# Row 66
$ws.Range("H66").Value = 3.6
$ws.Range("I66").Value = 1.8
$ws.Range("L66").Value = 1.24
$ws.Range("M66").Value = 3.65
$ws.Range("N66").Value = 1.72
$ws.Range("O66").Value = 2
$ws.Range("P66").Value = 1.36
$ws.Range("Q66").Value = 2.9
$ws.Range("S66").Value = 2.05
$ws.Range("T66").Value = 12.5
$ws.Range("X66").Value = 32
$ws.Range("Y66").Value = 35
$ws.Range("AB66").Value = 14
$ws.Range("AC66").Value = 55
$ws.Range("AD66").Value = 400
$ws.Range("AJ66").Value = 23
